# Generate Report for Handoff
#
# A new source file (bb1f29a8-f09c-4c46-868c-f33ab00afdc0.md) was picked up by the
# localization pipeline. This inserts a fresh "Ready for handoff" row for it right
# before the existing ".localization-config" / "Not to be localized" row on every
# sheet (Overview, zh-cn, de-de), pushing that row down by one, and relinks the
# hyperlinks for the affected rows.

$wb = $excel.ActiveWorkbook

$newFile        = "bb1f29a8-f09c-4c46-868c-f33ab00afdc0.md"
$repoCommit     = "29492cf99c376bf51e7b19b21261bb73e03c4323"
$handoffCommit  = "a5c06c06917c928260377249aa619c1044fc46d9"
$zhXlf          = "bb1f29a8-f09c-4c46-868c-f33ab00afdc0.$handoffCommit.zh-cn.xlf"
$deXlf          = "bb1f29a8-f09c-4c46-868c-f33ab00afdc0.$handoffCommit.de-de.xlf"
$zhDatetime     = "2016-03-01 05:58:43"
$deDatetime     = "2016-03-01 05:58:54"
$epoch          = "0001-01-01 00:00:00"

$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/e2e/$newFile"
$cfgUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/$repoCommit/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhXlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$handoffCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deXlf"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" (3 columns: File Name, zh-cn, de-de)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Push the old row 3 (".localization-config") down to row 4.
$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

# Write the new row 3 for the new source file.
$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# Rebuild the hyperlinks in row order.
$wsOverview.Cells.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", "a299f745-a960-40eb-9b6c-47ebe1fdc1d9.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdUrl.Replace($newFile, $newFile), "", "", $newFile) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Push the old row 3 (".localization-config" / "Ignored") down to row 4.
$wsZh.Range("A4").Value = ".localization-config"
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

# Write the new row 3 for the new source file.
$wsZh.Range("A3").Value = $newFile
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = $zhXlf
$wsZh.Range("D3").Value = $zhDatetime
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Include"

# Rebuild the hyperlinks in row order.
$wsZh.Cells.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", "a299f745-a960-40eb-9b6c-47ebe1fdc1d9.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5e5d06dcd82ea84b7c489eaf2b00ead4c417c081/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/a299f745-a960-40eb-9b6c-47ebe1fdc1d9.74f2cad9c7439f0f72e3a0640a9d2d86c5f00360.zh-cn.xlf", "", "", "a299f745-a960-40eb-9b6c-47ebe1fdc1d9.74f2cad9c7439f0f72e3a0640a9d2d86c5f00360.zh-cn.xlf") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl, "", "", $newFile) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zhXlfUrl, "", "", $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Push the old row 3 (".localization-config" / "Ignored") down to row 4.
$wsDe.Range("A4").Value = ".localization-config"
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

# Write the new row 3 for the new source file.
$wsDe.Range("A3").Value = $newFile
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = $deXlf
$wsDe.Range("D3").Value = $deDatetime
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Include"

# Rebuild the hyperlinks in row order.
$wsDe.Cells.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", "a299f745-a960-40eb-9b6c-47ebe1fdc1d9.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/eb1021545cc17703dac6b937c732eb3afa54ade6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/a299f745-a960-40eb-9b6c-47ebe1fdc1d9.74f2cad9c7439f0f72e3a0640a9d2d86c5f00360.de-de.xlf", "", "", "a299f745-a960-40eb-9b6c-47ebe1fdc1d9.74f2cad9c7439f0f72e3a0640a9d2d86c5f00360.de-de.xlf") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl, "", "", $newFile) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $deXlfUrl, "", "", $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $cfgUrl, "", "", ".localization-config") | Out-Null

Write-Host "Handoff report rows generated for" $newFile
